$d = $word.ActiveDocument

# Locate the target run's text and replace with the new text.
# Original: ", Trello, GitHub, y un entorno de desarrollo web basado en Node.js, Express y MySQL."
# New:      ", GitHub, y un entorno de desarrollo web basado en Python, Django y SQLite(en uso de pruebas aun)."

$found = $d.Content.Find.Execute(
    ", Trello, GitHub, y un entorno de desarrollo web basado en Node.js, Express y MySQL.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", GitHub, y un entorno de desarrollo web basado en Python, Django y SQLite(en uso de pruebas aun).",
    2)

if (-not $found) {
    throw "Target sentence not found; document may already be modified."
}
